# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (new DAMSLTag, new DialogAct)
$updates = @{
    10 = @("sd", "Statement-non-opinion")
    11 = @("sv", "Statement-opinion")
    12 = @("%", "Uninterpretable")
    14 = @("sd", "Statement-non-opinion")
    19 = @("sd", "Statement-non-opinion")
    21 = @("sv", "Statement-opinion")
    52 = @("sd", "Statement-non-opinion")
    57 = @("sd", "Statement-non-opinion")
    60 = @("sv", "Statement-opinion")
    70 = @("aa", "Agree/Accept")
    84 = @("sd", "Statement-non-opinion")
    91 = @("sv", "Statement-opinion")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
